$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column B (Notes column) entirely - shifts nothing left since it's the last column
$ws.Columns.Item(2).Delete()

# Rename header from "#SampleID" to "Sample_ID"
$ws.Range("A1").Value = "Sample_ID"

# Set the active cell/selection to A2 as in the after state
$ws.Range("A2").Select()
